$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Swap the Execute values: Guru99 test case -> No, OrangeHRM test case -> Yes
$ws.Range("B2").Value = "No"
$ws.Range("B3").Value = "Yes"

# Update the active selection on the sheet
$ws.Range("D4").Select()

# Update the workbook window position/size
$excel.ActiveWindow.WindowState = -4143
$excel.Windows.Item(1).Top = 330
$excel.Windows.Item(1).Left = 0
$excel.Windows.Item(1).Width = 14415
$excel.Windows.Item(1).Height = 2415
